$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 39 appended to the bitcoin_buys log (run on 2025-08-10).
# Force the date column to literal text so it isn't auto-converted to a
# date serial by Excel's input parser, then clear the temporary format so
# no stray style index is left behind on the cell.
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "08/10/2025"
$ws.Range("A39").ClearFormats()

$ws.Range("B39").Value = 0.0004207599999999992
$ws.Range("C39").Value = 118832.5886491114
$ws.Range("D39").Value = 50
